$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A4: was blank with a (now unused) date-format style; give it the value 9
# and a plain General format with a yellow highlight fill.
$ws.Range("A4").ClearFormats()
$ws.Range("A4").Value2 = 9
$ws.Range("A4").Interior.Color = 65535

# B4: give it the same text as B2 ("=>blank"); the leading apostrophe
# forces text entry (same quote-prefixed style Excel applies to B2).
$ws.Range("B4").Value2 = "'=>blank"

# Move the active selection to B4, matching the edited workbook.
$ws.Range("B4").Select()
